{"js": "// Update the multiplication-table answer cells in place.\n// Each entry is a unique \"A\u00d7B=C\" string that appears exactly once in the\n// document body, so a plain search-and-replace (no wildcards) is safe.\nconst replacements = [\n  [\"53\u00d786=4558\", \"32\u00d736=1152\"],\n  [\"45\u00d779=3555\", \"11\u00d758=638\"],\n  [\"31\u00d786=2666\", \"32\u00d767=2144\"],\n  [\"60\u00d730=1800\", \"95\u00d746=4370\"],\n  [\"70\u00d785=5950\", \"70\u00d790=6300\"],\n  [\"21\u00d765=1365\", \"74\u00d734=2516\"],\n  [\"19\u00d743=817\", \"52\u00d762=3224\"],\n  [\"71\u00d720=1420\", \"11\u00d756=616\"],\n  [\"35\u00d767=2345\", \"31\u00d753=1643\"],\n  [\"27\u00d791=2457\", \"86\u00d744=3784\"],\n  [\"70\u00d757=3990\", \"21\u00d718=378\"],\n  [\"83\u00d716=1328\", \"73\u00d762=4526\"],\n  [\"59\u00d729=1711\", \"59\u00d720=1180\"],\n  [\"58\u00d736=2088\", \"89\u00d763=5607\"],\n  [\"85\u00d779=6715\", \"21\u00d760=1260\"],\n  [\"93\u00d722=2046\", \"43\u00d718=774\"],\n  [\"48\u00d735=1680\", \"74\u00d768=5032\"],\n  [\"36\u00d796=3456\", \"44\u00d758=2552\"],\n  [\"73\u00d746=3358\", \"95\u00d767=6365\"],\n  [\"25\u00d757=1425\", \"86\u00d755=4730\"],\n  [\"31\u00d719=589\", \"32\u00d773=2336\"],\n  [\"46\u00d749=2254\", \"83\u00d779=6557\"],\n  [\"22\u00d712=264\", \"39\u00d765=2535\"],\n  [\"89\u00d755=4895\", \"37\u00d791=3367\"],\n  [\"12\u00d772=864\", \"22\u00d744=968\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the multiplication-table answer cells in place.\n# Each entry is a unique \"A\u00d7B=C\" string that appears exactly once in the\n# document body, so a plain Find/Replace (no wildcards) is safe for each.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"53\u00d786=4558\", \"32\u00d736=1152\"),\n    @(\"45\u00d779=3555\", \"11\u00d758=638\"),\n    @(\"31\u00d786=2666\", \"32\u00d767=2144\"),\n    @(\"60\u00d730=1800\", \"95\u00d746=4370\"),\n    @(\"70\u00d785=5950\", \"70\u00d790=6300\"),\n    @(\"21\u00d765=1365\", \"74\u00d734=2516\"),\n    @(\"19\u00d743=817\", \"52\u00d762=3224\"),\n    @(\"71\u00d720=1420\", \"11\u00d756=616\"),\n    @(\"35\u00d767=2345\", \"31\u00d753=1643\"),\n    @(\"27\u00d791=2457\", \"86\u00d744=3784\"),\n    @(\"70\u00d757=3990\", \"21\u00d718=378\"),\n    @(\"83\u00d716=1328\", \"73\u00d762=4526\"),\n    @(\"59\u00d729=1711\", \"59\u00d720=1180\"),\n    @(\"58\u00d736=2088\", \"89\u00d763=5607\"),\n    @(\"85\u00d779=6715\", \"21\u00d760=1260\"),\n    @(\"93\u00d722=2046\", \"43\u00d718=774\"),\n    @(\"48\u00d735=1680\", \"74\u00d768=5032\"),\n    @(\"36\u00d796=3456\", \"44\u00d758=2552\"),\n    @(\"73\u00d746=3358\", \"95\u00d767=6365\"),\n    @(\"25\u00d757=1425\", \"86\u00d755=4730\"),\n    @(\"31\u00d719=589\", \"32\u00d773=2336\"),\n    @(\"46\u00d749=2254\", \"83\u00d779=6557\"),\n    @(\"22\u00d712=264\", \"39\u00d765=2535\"),\n    @(\"89\u00d755=4895\", \"37\u00d791=3367\"),\n    @(\"12\u00d772=864\", \"22\u00d744=968\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $r = $d.Content\n    $r.Find.ClearFormatting()\n    $r.Find.Replacement.ClearFormatting()\n    $r.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
